$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
[void]$ws1.Range("D8:D201").Validation.Delete()
$ws1.Range("D8:D92").Validation.Add(3, 1, 1, "=Hidden_13")
$ws1.Range("D8:D92").Validation.IgnoreBlank = $true
$ws1.Range("D8:D92").Validation.ShowInput = $false
$ws1.Range("D8:D92").Validation.ShowError = $true
